# Update LinkedIn carousel slide copy from "renewable energy sector" generic
# messaging to the NTPC/GAIL joint-venture announcement.

$p = $ppt.ActivePresentation

$updates = @(
    @{
        Slide = 1
        Paragraphs = @(
            "NTPC and GAIL Joint Venture",
            "NTPC Green Energy has established a 50:50 joint venture (JV) with GAIL.",
            "The JV has received the necessary approvals."
        )
    },
    @{
        Slide = 2
        Paragraphs = @(
            "Investment Focus",
            "The partnership aims to enhance green energy production.",
            "Both companies are focusing on sustainable energy solutions."
        )
    },
    @{
        Slide = 3
        Paragraphs = @(
            "Strategic Goals",
            "The collaboration will leverage NTPC's renewable energy expertise and GAIL's infrastructure.",
            "The JV is expected to drive growth in the green energy sector."
        )
    },
    @{
        Slide = 4
        Paragraphs = @(
            "Market Impact",
            "NTPC Green Energy shares have gained attention following the JV announcement.",
            "Investor interest in green energy sectors is likely to increase."
        )
    },
    @{
        Slide = 5
        Paragraphs = @(
            "Future Prospects",
            "This JV aligns with India's commitment to expanding its renewable energy portfolio.",
            "Both companies plan to explore various green energy initiatives."
        )
    },
    @{
        Slide = 6
        Paragraphs = @(
            "Collaborative Benefits",
            "The joint ventures aim to capitalize on each company's strengths.",
            "Combining resources is expected to improve project efficiencies."
        )
    }
)

foreach ($u in $updates) {
    $s = $p.Slides.Item($u.Slide)
    $sh = $s.Shapes.Item(1)
    $tr = $sh.TextFrame.TextRange
    for ($i = 0; $i -lt $u.Paragraphs.Count; $i++) {
        $tr.Paragraphs($i + 1).Runs(1).Text = $u.Paragraphs[$i]
    }
}
